$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata row (row 8) to the new timestamp.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Set "Case Sensitive" (row 20) value to the literal text "true" (not a
# boolean). A leading apostrophe forces Excel to store it as text rather
# than auto-converting it to the Boolean TRUE value.
$ws.Range("B20").Value = "'true"
